# Update the cached "today" text shown by the Date placeholder's
# datetimeFigureOut field from 11/6/2022 to 11/7/2022.
#
# This field lives once on the Slide Master and once on each of the
# five Slide Layouts (PowerPoint keeps an independently-cached copy of
# the auto date field's display text per master/layout). We update the
# Date placeholder shape's text on the master and on every layout so
# the whole deck is consistent, mirroring what PowerPoint itself does
# when it recomputes the "Update automatically" date on save.

$p = $ppt.ActivePresentation

$oldDate = "11/6/2022"
$newDate = "11/7/2022"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
